$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.596.48"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "3.476.16"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'593.29"
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("D6").Value = "'167.73"
$ws.Range("E6").Value = "  -2.45%  "
$ws.Range("D7").Value = "'0.607"
$ws.Range("E7").Value = "  -2.34%  "
$ws.Range("D8").Value = "3.468.79"
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'0.192"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").Value = "'6.83"
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").Value = "'0.572"
$ws.Range("E12").Value = "  -5.24%  "
$ws.Range("D13").Value = "'46.57"
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "4.028.10"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("E16").Value = "  -5.98%  "
$ws.Range("D17").Value = "'612.65"
$ws.Range("E17").Value = "  -10.88%  "
$ws.Range("D18").Value = "3.477.24"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "68.627.81"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").Value = "'11.09"
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("E23").Value = "  -4.25%  "
$ws.Range("D24").Value = "'15.75"
$ws.Range("E24").Value = "  -5.18%  "
$ws.Range("D25").Value = "'95.67"
$ws.Range("E25").Value = "  -2.27%  "
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("D27").Value = "'5.83"
$ws.Range("E27").Value = "  +1.60%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("D30").Value = "'9.08"
$ws.Range("E30").Value = "  -3.78%  "
$ws.Range("D31").Value = "'32.80"
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("D32").Value = "'8.39"
$ws.Range("E32").Value = "  -5.29%  "
$ws.Range("E33").Value = "  -3.56%  "
$ws.Range("D35").Value = "'6.78"
$ws.Range("E35").Value = "  -6.52%  "
$ws.Range("D36").Value = "'570.75"
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("D37").Value = "'10.67"
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("E38").Value = "  -5.34%  "
$ws.Range("D39").Value = "'56.89"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("E40").Value = "  -4.42%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("D43").Value = "'0.0436"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").Value = "3.388.56"
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("E45").Value = "  -4.78%  "
$ws.Range("D46").Value = "'32.45"
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("D47").Value = "0.0₃0693"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("D48").Value = "'2.82"
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("E49").Value = "  -1.92%  "
$ws.Range("E50").Value = "  -4.31%  "
$ws.Range("D51").Value = "'132.94"
$ws.Range("E51").Value = "  -1.03%  "
